# Fonds de solidarite - Volet 2 - add 2020-09-30 data
# Updates nombre_aides (C) and montant_total (D) for the affected region/categorie-juridique rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("C2:D2").NumberFormat = "@"
$ws.Range("C2").Value = "1457"
$ws.Range("D2").Value = "3575243.70"

# Row 4
$ws.Range("C4:D4").NumberFormat = "@"
$ws.Range("C4").Value = "1116"
$ws.Range("D4").Value = "4465207.81"

# Row 6
$ws.Range("C6:D6").NumberFormat = "@"
$ws.Range("C6").Value = "757"
$ws.Range("D6").Value = "2724526.29"

# Row 7
$ws.Range("C7:D7").NumberFormat = "@"
$ws.Range("C7").Value = "20"
$ws.Range("D7").Value = "44998.41"

# Row 8
$ws.Range("C8:D8").NumberFormat = "@"
$ws.Range("C8").Value = "40"
$ws.Range("D8").Value = "180144.45"

# Row 14
$ws.Range("C14:D14").NumberFormat = "@"
$ws.Range("C14").Value = "249"
$ws.Range("D14").Value = "667219.00"

# Row 16
$ws.Range("C16:D16").NumberFormat = "@"
$ws.Range("C16").Value = "529"
$ws.Range("D16").Value = "2033649.55"

# Row 20
$ws.Range("C20:D20").NumberFormat = "@"
$ws.Range("C20").Value = "205"
$ws.Range("D20").Value = "551394.00"

# Row 22
$ws.Range("C22:D22").NumberFormat = "@"
$ws.Range("C22").Value = "365"
$ws.Range("D22").Value = "1562254.69"

# Row 23
$ws.Range("C23:D23").NumberFormat = "@"
$ws.Range("C23").Value = "172"
$ws.Range("D23").Value = "635289.45"

# Row 30
$ws.Range("C30:D30").NumberFormat = "@"
$ws.Range("C30").Value = "345"
$ws.Range("D30").Value = "909441.89"

# Row 32
$ws.Range("C32:D32").NumberFormat = "@"
$ws.Range("C32").Value = "664"
$ws.Range("D32").Value = "3282831.87"

# Row 34
$ws.Range("C34:D34").NumberFormat = "@"
$ws.Range("C34").Value = "438"
$ws.Range("D34").Value = "1664744.99"

# Row 42
$ws.Range("C42:D42").NumberFormat = "@"
$ws.Range("C42").Value = "275"
$ws.Range("D42").Value = "800923.15"

# Row 43
$ws.Range("C43:D43").NumberFormat = "@"
$ws.Range("C43").Value = "128"
$ws.Range("D43").Value = "731049.98"

# Row 44
$ws.Range("C44:D44").NumberFormat = "@"
$ws.Range("C44").Value = "190"
$ws.Range("D44").Value = "833553.25"

# Row 46
$ws.Range("C46:D46").NumberFormat = "@"
$ws.Range("C46").Value = "12"
$ws.Range("D46").Value = "70085.23"

# Row 47
$ws.Range("C47:D47").NumberFormat = "@"
$ws.Range("C47").Value = "539"
$ws.Range("D47").Value = "1519521.23"

# Row 49
$ws.Range("C49:D49").NumberFormat = "@"
$ws.Range("C49").Value = "763"
$ws.Range("D49").Value = "3591417.57"

# Row 50
$ws.Range("C50:D50").NumberFormat = "@"
$ws.Range("C50").Value = "536"
$ws.Range("D50").Value = "2227039.97"

# Row 51
$ws.Range("C51:D51").NumberFormat = "@"
$ws.Range("C51").Value = "7"
$ws.Range("D51").Value = "17760.00"

# Row 52
$ws.Range("C52:D52").NumberFormat = "@"
$ws.Range("C52").Value = "29"
$ws.Range("D52").Value = "163011.07"

# Row 53
$ws.Range("C53:D53").NumberFormat = "@"
$ws.Range("C53").Value = "6401"
$ws.Range("D53").Value = "15685066.28"

# Row 57
$ws.Range("C57:D57").NumberFormat = "@"
$ws.Range("C57").Value = "4875"
$ws.Range("D57").Value = "19322073.30"

# Row 70
$ws.Range("C70:D70").NumberFormat = "@"
$ws.Range("C70").Value = "56"
$ws.Range("D70").Value = "199560.00"

# Row 71
$ws.Range("C71:D71").NumberFormat = "@"
$ws.Range("C71").Value = "26"
$ws.Range("D71").Value = "109496.39"

# Row 73
$ws.Range("C73:D73").NumberFormat = "@"
$ws.Range("C73").Value = "262"
$ws.Range("D73").Value = "757179.00"

# Row 74
$ws.Range("C74:D74").NumberFormat = "@"
$ws.Range("C74").Value = "423"
$ws.Range("D74").Value = "1636446.36"

# Row 75
$ws.Range("C75:D75").NumberFormat = "@"
$ws.Range("C75").Value = "249"
$ws.Range("D75").Value = "1031641.51"

# Row 77
$ws.Range("C77:D77").NumberFormat = "@"
$ws.Range("C77").Value = "20"
$ws.Range("D77").Value = "81105.00"

# Row 78
$ws.Range("C78:D78").NumberFormat = "@"
$ws.Range("C78").Value = "420"
$ws.Range("D78").Value = "1112329.60"

# Row 80
$ws.Range("C80:D80").NumberFormat = "@"
$ws.Range("C80").Value = "1011"
$ws.Range("D80").Value = "4243442.14"

# Row 81
$ws.Range("C81:D81").NumberFormat = "@"
$ws.Range("C81").Value = "559"
$ws.Range("D81").Value = "2253293.12"

# Row 82
$ws.Range("C82:D82").NumberFormat = "@"
$ws.Range("C82").Value = "43"
$ws.Range("D82").Value = "142078.00"

# Row 84
$ws.Range("C84:D84").NumberFormat = "@"
$ws.Range("C84").Value = "730"
$ws.Range("D84").Value = "1837571.33"

# Row 87
$ws.Range("C87:D87").NumberFormat = "@"
$ws.Range("C87").Value = "1048"
$ws.Range("D87").Value = "4047515.84"

# Row 88
$ws.Range("C88:D88").NumberFormat = "@"
$ws.Range("C88").Value = "743"
$ws.Range("D88").Value = "2394406.92"

# Row 100
$ws.Range("C100:D100").NumberFormat = "@"
$ws.Range("C100").Value = "1253"
$ws.Range("D100").Value = "4579359.41"

# Row 102
$ws.Range("C102:D102").NumberFormat = "@"
$ws.Range("C102").Value = "1199"
$ws.Range("D102").Value = "4150973.98"
